$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list values to match the latest scrape.
# D/E columns are plain text cells (prices use '.' as thousands
# separators, e.g. '42.979.73', and percents keep their padding
# spaces) so we force text via NumberFormat '@' while writing,
# then ClearFormats() to drop back to the default (unstyled) cell
# so no stray per-cell formatting is introduced.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '42.979.73'
$c.ClearFormats()

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  +2.18%  '
$c.ClearFormats()

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.299.43'
$c.ClearFormats()

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  +1.59%  '
$c.ClearFormats()

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.ClearFormats()

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '301.45'
$c.ClearFormats()

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  +1.17%  '
$c.ClearFormats()

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  +5.23%  '
$c.ClearFormats()

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  +2.00%  '
$c.ClearFormats()

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c.ClearFormats()

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  +3.59%  '
$c.ClearFormats()

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '34.09'
$c.ClearFormats()

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  +3.58%  '
$c.ClearFormats()

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  +1.13%  '
$c.ClearFormats()

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '49.01'
$c.ClearFormats()

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  +1.38%  '
$c.ClearFormats()

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  +4.31%  '
$c.ClearFormats()

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '17.89'
$c.ClearFormats()

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  +15.35%  '
$c.ClearFormats()

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '6.78'
$c.ClearFormats()

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  +2.21%  '
$c.ClearFormats()

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '2.656.79'
$c.ClearFormats()

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  +1.64%  '
$c.ClearFormats()

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.259.22'
$c.ClearFormats()

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -0.84%  '
$c.ClearFormats()

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.808'
$c.ClearFormats()

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  +4.87%  '
$c.ClearFormats()

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '42.907.74'
$c.ClearFormats()

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  +2.03%  '
$c.ClearFormats()

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  +9.18%  '
$c.ClearFormats()

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  +1.96%  '
$c.ClearFormats()

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.08'
$c.ClearFormats()

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  +1.29%  '
$c.ClearFormats()

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '67.67'
$c.ClearFormats()

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  +1.65%  '
$c.ClearFormats()

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '236.33'
$c.ClearFormats()

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  +1.50%  '
$c.ClearFormats()

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  +14.92%  '
$c.ClearFormats()

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  +0.07%  '
$c.ClearFormats()

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  +0.75%  '
$c.ClearFormats()

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '24.69'
$c.ClearFormats()

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  +4.04%  '
$c.ClearFormats()

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '167.90'
$c.ClearFormats()

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  +0.43%  '
$c.ClearFormats()

$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -9.15%  '
$c.ClearFormats()

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '33.77'
$c.ClearFormats()

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  +0.80%  '
$c.ClearFormats()

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  +1.20%  '
$c.ClearFormats()

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  +0.07%  '
$c.ClearFormats()

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '5.03'
$c.ClearFormats()

$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  +2.49%  '
$c.ClearFormats()

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '4.55'
$c.ClearFormats()

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  +1.59%  '
$c.ClearFormats()

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.42'
$c.ClearFormats()

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  +3.45%  '
$c.ClearFormats()

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  +4.73%  '
$c.ClearFormats()

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  +0.32%  '
$c.ClearFormats()

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  +4.00%  '
$c.ClearFormats()

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  +4.90%  '
$c.ClearFormats()

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  +0.91%  '
$c.ClearFormats()

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  +0.36%  '
$c.ClearFormats()

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -4.18%  '
$c.ClearFormats()

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.989.12'
$c.ClearFormats()

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  +1.71%  '
$c.ClearFormats()

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  +2.53%  '
$c.ClearFormats()

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '9.99'
$c.ClearFormats()

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  +5.07%  '
$c.ClearFormats()

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '17.54'
$c.ClearFormats()

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  +1.82%  '
$c.ClearFormats()

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.84'
$c.ClearFormats()

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  +2.70%  '
$c.ClearFormats()

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '56.64'
$c.ClearFormats()

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  +9.82%  '
$c.ClearFormats()

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.524.78'
$c.ClearFormats()

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  +3.72%  '
$c.ClearFormats()

